# Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Pepino dulce
# Weekly data refresh: insert 4 new rows (date 44753) above the existing
# row 278 block, pushing the rest of the table down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 278 (old rows 278:295 shift down to 282:299).
$ws.Range("A278:A281").EntireRow.Insert()

# Common (constant) values for this market / category / variety block.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112043
$categoria = "Pepino dulce"
$variedad  = "Cultivar IV Región"
$unidad    = "`$/bandeja 18 kilos"
$origen    = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"
$fecha = 44753

# New row data: row -> (calidad, volumen, precioMin, precioMax, precioProm, precioKg)
$newRows = @(
    @{ Row = 278; Calidad = "Especial"; Volumen = 260; PMin = 15000; PMax = 15000; PProm = 15000; PKg = 833 },
    @{ Row = 279; Calidad = "Primera";  Volumen = 320; PMin = 13000; PMax = 13000; PProm = 13000; PKg = 722 },
    @{ Row = 280; Calidad = "Segunda";  Volumen = 230; PMin = 10000; PMax = 10000; PProm = 10000; PKg = 556 },
    @{ Row = 281; Calidad = "Tercera";  Volumen = 130; PMin = 8000;  PMax = 8000;  PProm = 8000;  PKg = 444 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $catId
    $ws.Cells.Item($r, 7).Value = $categoria
    $ws.Cells.Item($r, 8).Value = $variedad
    $ws.Cells.Item($r, 9).Value = $rowData.Calidad
    $ws.Cells.Item($r, 10).Value = $rowData.Volumen
    $ws.Cells.Item($r, 11).Value = $rowData.PMin
    $ws.Cells.Item($r, 12).Value = $rowData.PMax
    $ws.Cells.Item($r, 13).Value = $rowData.PProm
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $rowData.PKg
    $ws.Cells.Item($r, 17).Value = $kgUnidades
    $ws.Cells.Item($r, 18).Value = $clasificacion
}
